$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ball = "('Ball Lightning', ['{R}{R}{R}', 'Creature — Elemental', 'Trample (This creature can deal excess combat damage to the player or planeswalker it’s attacking.)', 'Haste (This creature can attack and {T} as soon as it comes under your control.)', 'At the beginning of the end step, sacrifice Ball Lightning.', '6/1'])"
$oath = "('Oath of Druids', ['{1}{G}', 'Enchantment', 'At the beginning of each player’s upkeep, that player chooses target player who controls more creatures than they do and is their opponent. The first player may reveal cards from the top of their library until they reveal a creature card. If the first player does, that player puts that card onto the battlefield and all other cards revealed this way into their graveyard.'])"

$ws.Range("A2").Value = $ball
$ws.Range("A3").Value = $oath

$ws.Range("A4:A12").EntireRow.Delete()
